# WhitePack_Full_Template_v7.2.xlsx — "Add files via upload"
#
# Semantic changes applied (per the OOXML diff):
#  1. Make "Products" the active sheet (previously "SpecialPrices" was active).
#  2. Products!J2 ("Image" column) is re-pointed from the placeholder
#     https://example.com/image.jpg to a real thumbnail URL, and the cell is
#     turned into a clickable hyperlink (adds the built-in "Hyperlink" cell
#     style / font to the workbook).
#  3. A few column widths are set on the Products sheet (B, C, J).
#  4. Selection/cursor bookmarks are updated: Products -> J6, SpecialPrices -> A3.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Products")
$ws3 = $wb.Worksheets.Item("SpecialPrices")

# --- 3. Column widths on the Products sheet ---------------------------------
$ws1.Columns.Item(2).ColumnWidth  = 12.5             # B
$ws1.Columns.Item(3).ColumnWidth  = 10.8333333333333  # C
$ws1.Columns.Item(10).ColumnWidth = 37.6666666666667  # J

# --- 2. Turn J2 into a hyperlink pointing at the new thumbnail URL ----------
$newUrl = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcQTRIz_5va2njvWVvsNsQyIYRvuqJTCAHoIqA&s"
$ws1.Hyperlinks.Add($ws1.Range("J2"), $newUrl, "", "", $newUrl) | Out-Null

# --- 4. Restore each sheet's remembered selection ---------------------------
$ws3.Range("A3").Select() | Out-Null

# --- 1. Activate "Products" (also clears tabSelected on SpecialPrices) -----
$ws1.Select() | Out-Null
$ws1.Range("J6").Select() | Out-Null
